$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1560082246840916
$ws.Range("C2").Value = 1.334666554170477
$ws.Range("D2").Value = 7.673042805661654
$ws.Range("E2").Value = 2.770025777075306
$ws.Range("F2").Value = 2.793148460838132
$ws.Range("G2").Value = 51
$ws.Range("B3").Value = 0.07663322156689945
$ws.Range("C3").Value = 1.299834471760893
$ws.Range("D3").Value = 7.386217108190144
$ws.Range("E3").Value = 2.7177595751262
$ws.Range("F3").Value = 2.744260142950399
$ws.Range("G3").Value = 50
$ws.Range("B4").Value = 0.1568606132576995
$ws.Range("C4").Value = 1.29112821977527
$ws.Range("D4").Value = 7.117402227092668
$ws.Range("E4").Value = 2.667845990137487
$ws.Range("F4").Value = 2.690829533573063
$ws.Range("G4").Value = 49
$ws.Range("B5").Value = 0.0935746678326199
$ws.Range("C5").Value = 1.365793217520896
$ws.Range("D5").Value = 7.57656095386177
$ws.Range("E5").Value = 2.752555349827097
$ws.Range("F5").Value = 2.78007587260341
$ws.Range("G5").Value = 48
$ws.Range("B6").Value = 0.1424348649439729
$ws.Range("C6").Value = 1.381738402566604
$ws.Range("D6").Value = 8.024261160333362
$ws.Range("E6").Value = 2.832712685807257
$ws.Range("F6").Value = 2.859715526557411
$ws.Range("G6").Value = 47
$ws.Range("B7").Value = 0.1066018923057719
$ws.Range("C7").Value = 1.373320158004856
$ws.Range("D7").Value = 7.768055166292034
$ws.Range("E7").Value = 2.787123098517903
$ws.Range("F7").Value = 2.815859037392981
$ws.Range("G7").Value = 46
$ws.Range("B8").Value = 0.1578286088130434
$ws.Range("C8").Value = 1.342716586736134
$ws.Range("D8").Value = 7.844528394684448
$ws.Range("E8").Value = 2.800808525173481
$ws.Range("F8").Value = 2.827956351813759
$ws.Range("G8").Value = 45
$ws.Range("B9").Value = 0.08991622227243078
$ws.Range("C9").Value = 1.431264513475095
$ws.Range("D9").Value = 8.13824447033133
$ws.Range("E9").Value = 2.852760850532573
$ws.Range("F9").Value = 2.884308066253469
$ws.Range("G9").Value = 44
$ws.Range("B10").Value = 0.157804923366503
$ws.Range("C10").Value = 1.414187432328667
$ws.Range("D10").Value = 8.224904023511264
$ws.Range("E10").Value = 2.867909347157135
$ws.Range("F10").Value = 2.897454014080703
$ws.Range("G10").Value = 43
$ws.Range("B11").Value = 0.06053529240133091
$ws.Range("C11").Value = 1.461861534950335
$ws.Range("D11").Value = 8.447652183675055
$ws.Range("E11").Value = 2.906484506009804
$ws.Range("F11").Value = 2.941077792346796
$ws.Range("G11").Value = 42
